$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.185.82"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.277.72"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "'114.39"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'304.91"
$ws.Range("E6").Value = "  +6.83%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.614"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'44.90"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "'55.09"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "'8.88"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("E14").Value = "  +18.83%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'15.41"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "2.621.87"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "2.284.77"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "43.159.89"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'7.23"
$ws.Range("E21").Value = "  +4.77%  "
$ws.Range("D22").Value = "'74.89"
$ws.Range("E22").Value = "  +2.56%  "
$ws.Range("D23").Value = "'3.55"
$ws.Range("E23").Value = "  +10.93%  "
$ws.Range("E24").Value = "  +4.66%  "
$ws.Range("D25").Value = "'255.09"
$ws.Range("E25").Value = "  +10.00%  "
$ws.Range("D26").Value = "'9.03"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").Value = "'11.72"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'38.20"
$ws.Range("E30").Value = "  -5.20%  "
$ws.Range("D31").Value = "'175.32"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'22.11"
$ws.Range("E32").Value = "  +4.38%  "
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("D34").Value = "'0.0898"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "'5.72"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("D36").Value = "'5.06"
$ws.Range("E36").Value = "  +9.14%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.129"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.27"
$ws.Range("E38").Value = "  -8.05%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "'2.53"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").Value = "'72.88"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'12.64"
$ws.Range("E45").Value = "  -6.39%  "
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").Value = "'5.63"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "'106.89"
$ws.Range("E48").Value = "  +5.57%  "
$ws.Range("D49").Value = "'8.79"
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'73.84"
$ws.Range("E51").Value = "  +5.59%  "
